$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.7474139358276786
$ws.Range("E2").Value = 0.7474139358276786

# Row 3
$ws.Range("D3").Value = 0.0003939682687483869
$ws.Range("E3").Value = 0.0003939682687483869

# Row 4
$ws.Range("D4").Value = 0.01559977197893658
$ws.Range("E4").Value = 0.01559977197893658

# Row 5
$ws.Range("D5").Value = 0.0002463296822664912
$ws.Range("E5").Value = 0.0002463296822664912

# Row 6
$ws.Range("D6").Value = 0.9583940106226873
$ws.Range("E6").Value = 0.9583940106226873

# Row 7
$ws.Range("D7").Value = 0.9924205679872576
$ws.Range("E7").Value = 0.00757943201274236

# Row 8
$ws.Range("D8").Value = 0.8703013505422813
$ws.Range("E8").Value = 0.1296986494577187

# Row 9
$ws.Range("D9").Value = 0.8017258803033372
$ws.Range("E9").Value = 0.1982741196966628

# Row 10
$ws.Range("D10").Value = 0.9868680309480609
$ws.Range("E10").Value = 0.01313196905193914

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.3875683154749018
$ws.Range("E11").Value = 0.6124316845250982
$ws.Range("F11").Value = 0.5900471806526184
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.9399071752257445
$ws.Range("E12").Value = 0.9399071752257445

# Row 13
$ws.Range("D13").Value = [double]"1.025645611479696E-05"
$ws.Range("E13").Value = [double]"1.025645611479696E-05"

# Row 14
$ws.Range("D14").Value = 0.003525057537844885
$ws.Range("E14").Value = 0.003525057537844885

# Row 15
$ws.Range("D15").Value = 0.0002182980614822648
$ws.Range("E15").Value = 0.0002182980614822648

# Row 16
$ws.Range("D16").Value = 0.9966281373826422
$ws.Range("E16").Value = 0.9966281373826422

# Row 17
$ws.Range("D17").Value = 0.9995934061220472
$ws.Range("E17").Value = 0.0004065938779528322

# Row 18
$ws.Range("D18").Value = 0.5680510187688258
$ws.Range("E18").Value = 0.4319489812311742

# Row 19
$ws.Range("D19").Value = 0.9979619122808111
$ws.Range("E19").Value = 0.002038087719188897

# Row 20
$ws.Range("D20").Value = 0.9995658757830147
$ws.Range("E20").Value = 0.0004341242169852988

# Row 21
$ws.Range("D21").Value = 0.990401166427946
$ws.Range("E21").Value = 0.009598833572054044
$ws.Range("F21").Value = 0.9085984230041504
